$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# The edit rotates several paragraph/run text values among each other
# (objectives <-> docentes list <-> programa resumido <-> avaliacao values,
# and separately programa <-> bibliografia). Doing a naive sequential
# Find/Replace would have later steps re-match text that an earlier step
# just inserted. To avoid that, first "tag" every source location with a
# unique placeholder token, then in a second pass replace each placeholder
# with its final destination text.

# --- Phase 1: tag each source location with a unique placeholder ---
Replace-Text "Fornecer os conhecimentos teóricos e práticos sobre estatística aplicada, materialografia e análises térmicas de materiais." "@@LOC_OBJETIVOS@@"
Replace-Text "519033 - Carlos Yujiro Shigue" "@@LOC_DOCENTE1@@"
Replace-Text "5840963 - Daniela Camargo Vernilli" "@@LOC_DOCENTE2@@"
Replace-Text "6495737 - Durval Rodrigues Junior" "@@LOC_DOCENTE3@@"
Replace-Text "984972 - Hugo Ricardo Zschommler Sandim" "@@LOC_DOCENTE4@@"
Replace-Text "Técnicas de Materialografia." "@@LOC_PROGRAMARESUMIDO@@"
Replace-Text "Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos." "@@LOC_METODO@@"
Replace-Text "Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0." "@@LOC_CRITERIO@@"
Replace-Text "Devido às características práticas da disciplina, não será oferecida recuperação." "@@LOC_NORMA@@"

# --- Phase 2: replace each placeholder with its final destination text ---
Replace-Text "@@LOC_OBJETIVOS@@" "Técnicas de Materialografia."
Replace-Text "@@LOC_DOCENTE1@@" "Fornecer os conhecimentos teóricos e práticos sobre estatística aplicada, materialografia e análises térmicas de materiais."
Replace-Text "@@LOC_DOCENTE2@@" "1. MATERIALOGRAFIA: Ensaio macrográfico ou macrografia; ensaio micrográfico ou micrografia. CORPO DE PROVA OU AMOSTRA. CORTE: discos de corte. Procedimento para o corte. EMBUTIMENTO: Preparação de corpo de prova: corpo de prova embutido a quente e a frio. Corpo de prova não embutido. LIXAMENTO: tipos de lixa; procedimento para o lixamento. POLIMENTO: processo mecânico; cuidados a serem observados no polimento. Processo semiautomático; processo eletrolítico; processo mecânico eletrolítico; polimento químico. Escolha do tipo de polimento. Procedimento para o polimento. ATAQUE QUÍMICO: princípio; métodos para obtenção de contraste. MICROSCOPIA ÓPTICA: Iluminação campo escuro; luz polarizada; contraste de fase; interferência diferencial. Partes de um microscópio óptico de reflexão; elementos mecânicos; elementos ópticos; iluminador; acessórios. Princípio da formação da imagem. Microscópio óptico de reflexão."
Replace-Text "@@LOC_DOCENTE3@@" "Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos."
Replace-Text "@@LOC_DOCENTE4@@" "Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."
Replace-Text "@@LOC_PROGRAMARESUMIDO@@" "Devido às características práticas da disciplina, não será oferecida recuperação."
Replace-Text "@@LOC_METODO@@" "519033 - Carlos Yujiro Shigue"
Replace-Text "@@LOC_CRITERIO@@" "5840963 - Daniela Camargo Vernilli"
Replace-Text "@@LOC_NORMA@@" "6495737 - Durval Rodrigues Junior"

# --- Structural swap: the "Programa" paragraph and the "Bibliografia"
#     paragraph trade their entire (single-run) content. Use Range.Text
#     with Chr(11) (vertical tab) for manual line breaks (<w:br/>). ---

# Locate the "Programa" heading paragraph and take the paragraph right after it
$paragraphs = $d.Paragraphs
$programaPara = $null
$bibliografiaPara = $null
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $txt = $paragraphs.Item($i).Range.Text
    if ($txt -eq "Programa`r") {
        $programaPara = $paragraphs.Item($i + 1)
    }
    if ($txt -eq "Bibliografia`r") {
        $bibliografiaPara = $paragraphs.Item($i + 1)
    }
}

$bibliografiaText = "COLPAERT; HUBERTUS. Metalografia dos produtos siderúrgicos comuns, 3ª Edição, Editora Edgard Blücher Ltda, São" + [char]11 + "Paulo – 1974." + [char]11 + "COUTINHO, TELMO DE AZEVEDO. Metalografia de Não-Ferrosos, Editora Edgard Blücher Ltda, São Paulo – 1980." + [char]11 + "PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985." + [char]11 + "MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001." + [char]11 + "WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008." + [char]11 + "REED-HILL, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982." + [char]11 + "Nondestructive Characterization of Materials. Series. Plenum Press, New York." + [char]11 + "YACOBI, B.G.; HOLT, D.B.; KAZMERSKI, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994."

$hugoText = "984972 - Hugo Ricardo Zschommler Sandim"

$programaPara.Range.Text = $bibliografiaText
$bibliografiaPara.Range.Text = $hugoText
